$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fermentation improvements in scenarios C and D:
# Fermentation TAL yield (row 19): E19 0.68 -> 0.73
$ws.Range("E19").Value = 0.73

# Fermentation TAL titer (row 20): E20 76 -> 68
$ws.Range("E20").Value = 68

# Update the current selection to match the saved workbook state (K19)
$ws.Activate()
$ws.Range("K19").Select()
